$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SwateTemplateMetadata")

# Rename the metadata sheet
$ws.Name = "isa_template"

# Clear the cell-formatting (border/fill) that used to frame the
# Tags / Tags Term Accession Number / Tags Term Source REF block,
# while keeping the cell contents intact.
$ws.Range("C12:E12").ClearFormats()
$ws.Range("C13").ClearFormats()

# These cells were empty placeholders that only carried formatting -
# clearing removes them completely (no value, no style).
$ws.Range("E13").Clear()
$ws.Range("C14:D14").Clear()
